# Insert a new data row at row 163, shifting existing rows 163:269 down to 164:270,
# then populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 163..269 down by inserting a new blank row at 163.
$ws.Rows.Item(163).Insert()

# Populate the new row 163 with the new "Primera" record (dated 2021-10-19 = serial 44488).
$ws.Cells.Item(163, 1).Value = 7
$ws.Cells.Item(163, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(163, 3).Value = "Ñuble"
$ws.Cells.Item(163, 4).Value = 44488
$ws.Cells.Item(163, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(163, 5).Value = 16
$ws.Cells.Item(163, 6).Value = 100114014
$ws.Cells.Item(163, 7).Value = "Betarraga"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 200
$ws.Cells.Item(163, 11).Value = 750
$ws.Cells.Item(163, 12).Value = 800
$ws.Cells.Item(163, 13).Value = 775
$ws.Cells.Item(163, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 155
$ws.Cells.Item(163, 17).Value = 5
$ws.Cells.Item(163, 18).Value = "Hortaliza"
